$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.957.84"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.71"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.92"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4660"
$ws.Range("E7").Value = "  +0.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  -0.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07360"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8725"
$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.27"
$ws.Range("E11").Value = "  -0.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.809.26"
$ws.Range("E12").Value = "  -0.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.415"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07107"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.522"
$ws.Range("E15").Value = "  +0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.55"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008706"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.67"
$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.983.22"

$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.079.08"
$ws.Range("E24").Value = "  +3.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.893"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.86"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.144"
$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.279"
$ws.Range("E29").Value = "  -0.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.04"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08873"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7588"
$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.506"
$ws.Range("E34").Value = "  +1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.909"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.096"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05300"
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("E39").Value = "  -0.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.973"
$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.186"
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5296"
$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.331"
$ws.Range("E43").Value = "  -4.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1656"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.447"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4872"

$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("E48").Value = "  +0.15%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.53"
$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.664"
$ws.Range("E50").Value = "  -0.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06294"
$ws.Range("E51").Value = "  +0.06%  "
